$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
#    This shared string is used as the "Status" value for every row on
#    all three sheets (Overview uses it per-language in columns E/F,
#    while the per-language sheets use it in column C).
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"

# ----------------------------------------------------------------------
# 2) Narrow the "Status" columns (now that the text is shorter) from
#    ~17.22 characters wide down to ~13.41 characters wide.
# ----------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
